$d = $word.ActiveDocument

# 1) Apply strike-through formatting to the "Numbers in Excel output is in 3 digits"
#    bullet paragraph (paragraph mark + every run get <w:strike/>).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Numbers in Excel output*in 3 digits*") {
        $p.Range.Font.StrikeThrough = 1
    }
}

# 2) Remove the two "# Load the required packages" / "library(tidyverse)" paragraphs
#    entirely (including their paragraph marks), leaving the surrounding blank
#    paragraphs untouched.
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*# Load the required packages*") {
        $startPara = $p
    }
    if ($t -like "*library(*tidyverse*)*") {
        $endPara = $p
    }
}

if (($startPara -ne $null) -and ($endPara -ne $null)) {
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
